# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.146.72"
$ws.Range("E2").Value = "  -5.55%  "
$ws.Range("D3").Value = "2.997.64"
$ws.Range("E3").Value = "  -5.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'572.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").Value = "'125.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.58%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "2.991.97"
$ws.Range("E8").Value = "  -6.09%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -8.49%  "
$ws.Range("E11").Value = "  -6.25%  "
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("E13").Value = "  -8.99%  "
$ws.Range("D14").Value = "'32.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.38%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "3.491.22"
$ws.Range("E16").Value = "  -5.89%  "
$ws.Range("D17").Value = "2.992.60"
$ws.Range("E17").Value = "  -6.04%  "
$ws.Range("D18").Value = "60.087.49"
$ws.Range("E18").Value = "  -5.60%  "
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "'427.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.72%  "
$ws.Range("E21").Value = "  -6.56%  "
$ws.Range("D22").Value = "'0.667"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.05%  "
$ws.Range("E23").Value = "  -8.52%  "
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "'79.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.06%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -5.61%  "
$ws.Range("E29").Value = "  -7.18%  "
$ws.Range("E30").Value = "  -7.27%  "
$ws.Range("D31").Value = "'6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.29%  "
$ws.Range("D32").Value = "'25.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.59%  "
$ws.Range("D33").Value = "'0.0940"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("D35").Value = "'0.936"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'50.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.10%  "
$ws.Range("D38").Value = "0.0₃0666"
$ws.Range("E38").Value = "  -9.95%  "
$ws.Range("D39").Value = "'8.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("E40").Value = "  -9.99%  "
$ws.Range("E41").Value = "  -5.92%  "
$ws.Range("D42").Value = "'375.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("D43").Value = "2.677.82"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("E44").Value = "  -8.21%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -7.89%  "
$ws.Range("D47").Value = "'120.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.19%  "
$ws.Range("D48").Value = "'2.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.16%  "
$ws.Range("E49").Value = "  -4.05%  "
$ws.Range("D50").Value = "'23.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.49%  "
$ws.Range("E51").Value = "  -2.19%  "
